# Araklı 1961 Spor vs Çirihtalar - Maç Sonucu eklendi
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Playoff")

# Match result: Çirihtalar (A9) 1 - 5 Araklı 1961 Spor (B9)
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 5

# Move the active cell selection to E13
$ws.Range("E13").Select()
